# Scheduled-runner refresh of leve-profit figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the job sheets, as pulled by the
# Kraken_Profits market-data sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 95.25
$ws.Range("I2").Value = 103
$ws.Range("J2").Value = 87.5
$ws.Range("K2").Value = 103
$ws.Range("L2").Value = 87.5
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = -313.5
$ws.Range("H4").Value = 2560
$ws.Range("I4").Value = 2560
$ws.Range("K4").Value = 2560
$ws.Range("M4").Value = -2446
$ws.Range("H5").Value = 104.94444
$ws.Range("I5").Value = 90
$ws.Range("K5").Value = 90
$ws.Range("M5").Value = 25
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("M12").Value = 70
$ws.Range("H70").Value = 2863.6365
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 3375
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 10125
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -10665
$ws.Range("H73").Value = 2863.6365
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 3375
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 10125
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -11997
$ws.Range("H113").Value = 1651
$ws.Range("I113").Value = 1651
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1651
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1603
$ws.Range("N113").ClearContents()
$ws.Range("H138").Value = 6219.1665
$ws.Range("J138").Value = 5388
$ws.Range("L138").Value = 16164
$ws.Range("N138").Value = -26444

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13900.667
$ws.Range("I32").Value = 11684.546
$ws.Range("K32").Value = 11684.546
$ws.Range("M32").Value = -11397.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 239.54546
$ws.Range("I7").Value = 289
$ws.Range("J7").Value = 198.33333
$ws.Range("K7").Value = 289
$ws.Range("L7").Value = 198.33333
$ws.Range("M7").Value = -176
$ws.Range("N7").Value = -424.33333
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1480
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1384
$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877
$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 143228.14
$ws.Range("I4").Value = 319
$ws.Range("J4").Value = 200391.8
$ws.Range("K4").Value = 957
$ws.Range("L4").Value = 601175.3999999999
$ws.Range("M4").Value = -845
$ws.Range("N4").Value = -601399.3999999999
$ws.Range("H7").Value = 16666838
$ws.Range("I7").Value = 3.5
$ws.Range("K7").Value = 10.5
$ws.Range("M7").Value = 101.5
$ws.Range("H74").Value = 969
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 969
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 2907
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -5029
$ws.Range("H77").Value = 969
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 969
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 8721
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -19329
$ws.Range("H122").Value = 5821.5557
$ws.Range("J122").Value = 5033.8
$ws.Range("L122").Value = 45304.2
$ws.Range("N122").Value = -50204.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2780200.5
$ws.Range("I14").Value = 5001334.5
$ws.Range("J14").Value = 1669633.4
$ws.Range("K14").Value = 5001334.5
$ws.Range("L14").Value = 1669633.4
$ws.Range("M14").Value = -5001166.5
$ws.Range("N14").Value = -1669969.4
$ws.Range("H70").Value = 3994.5
$ws.Range("I70").Value = 3994.5
$ws.Range("K70").Value = 3994.5
$ws.Range("M70").Value = -3724.5
$ws.Range("H73").Value = 3994.5
$ws.Range("I73").Value = 3994.5
$ws.Range("K73").Value = 3994.5
$ws.Range("M73").Value = -3058.5
$ws.Range("H80").Value = 2749.8572
$ws.Range("I80").Value = 2569.8
$ws.Range("K80").Value = 2569.8
$ws.Range("M80").Value = -1571.8
$ws.Range("H83").Value = 2749.8572
$ws.Range("I83").Value = 2569.8
$ws.Range("K83").Value = 12849
$ws.Range("M83").Value = -7857
$ws.Range("H99").Value = 11887.5
$ws.Range("I99").Value = 5025
$ws.Range("K99").Value = 5025
$ws.Range("M99").Value = -2779

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 741.1111
$ws.Range("I22").Value = 678.3333
$ws.Range("J22").Value = 866.6667
$ws.Range("K22").Value = 678.3333
$ws.Range("L22").Value = 866.6667
$ws.Range("M22").Value = -383.3333
$ws.Range("N22").Value = -1456.6667
$ws.Range("H27").Value = 741.1111
$ws.Range("I27").Value = 678.3333
$ws.Range("J27").Value = 866.6667
$ws.Range("K27").Value = 678.3333
$ws.Range("L27").Value = 866.6667
$ws.Range("M27").Value = -571.3333
$ws.Range("N27").Value = -1080.6667
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H46").Value = 840
$ws.Range("I46").Value = 733.3333
$ws.Range("K46").Value = 733.3333
$ws.Range("M46").Value = -545.3333
$ws.Range("H82").Value = 1612.1666
$ws.Range("I82").Value = 1780.875
$ws.Range("J82").Value = 1274.75
$ws.Range("K82").Value = 1780.875
$ws.Range("L82").Value = 1274.75
$ws.Range("M82").Value = -1419.875
$ws.Range("N82").Value = -1996.75
$ws.Range("H85").Value = 1612.1666
$ws.Range("I85").Value = 1780.875
$ws.Range("J85").Value = 1274.75
$ws.Range("K85").Value = 1780.875
$ws.Range("L85").Value = 1274.75
$ws.Range("M85").Value = -532.875
$ws.Range("N85").Value = -3770.75
$ws.Range("H136").Value = 8862.909
$ws.Range("J136").Value = 19833.334
$ws.Range("L136").Value = 59500.00199999999
$ws.Range("N136").Value = -64600.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39999
$ws.Range("I70").Value = 39999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 39999
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -39684
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 39999
$ws.Range("I73").Value = 39999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 39999
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -38907
$ws.Range("N73").ClearContents()
$ws.Range("H94").Value = 26499.5
$ws.Range("I94").Value = 19999
$ws.Range("K94").Value = 19999
$ws.Range("M94").Value = -19098
$ws.Range("H136").Value = 2224.8333
$ws.Range("I136").Value = 2224.8333
$ws.Range("K136").Value = 6674.499899999999
$ws.Range("M136").Value = -4124.499899999999
